# Row 7 and row 8 in the "Artfynd" sheet represent two separate tree-branch
# sighting records that were merged/duplicated with the wrong Id assigned to
# the wrong set of location/observer/substrate details. This script swaps the
# per-record fields between row 7 and row 8 so each record carries its
# correct data, leaving the fields that are identical between the two rows
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture current ("before") values for the columns that differ between
#     row 7 and row 8 ---
$cols = @("A","J","K","N","P","Q","R","AC","AF","AH","AJ","AK","AM","AO","AW","AX")

$row7 = @{}
$row8 = @{}
foreach ($col in $cols) {
    $row7[$col] = $ws.Range("$col`7").Value2
    $row8[$col] = $ws.Range("$col`8").Value2
}

# --- write row 8's original values into row 7, and vice versa ---
foreach ($col in $cols) {
    $ws.Range("$col`7").Value2 = $row8[$col]
    $ws.Range("$col`8").Value2 = $row7[$col]
}
